$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -5.981199999999998
$ws.Range("D8").Value = -8.84059999999999
$ws.Range("C12").Value = -14.75070000000001
$ws.Range("D12").Value = -8.152000000000001
$ws.Range("D14").Value = -8.735799999999998
$ws.Range("D22").Value = -7.900899999999996
